$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cfh"
$ws.Cells.Item(2, 3).Value = "Sell"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.362097333333333
$ws.Cells.Item(2, 8).Value = 4.086292
$ws.Cells.Item(2, 9).Value = 0.0196292037450248
$ws.Cells.Item(2, 10).Value = 0.0196292037450248
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2056386666666667
$ws.Cells.Item(2, 14).Value = 0.616916
$ws.Cells.Item(2, 15).Value = 0.004186411275012692
$ws.Cells.Item(2, 16).Value = 0.004186411275012692
$ws.Cells.Item(2, 17).Value = 0.2800998794968889
$ws.Cells.Item(2, 18).Value = 2.520898915472
$ws.Cells.Item(2, 19).Value = 0.00008217591987769319
$ws.Cells.Item(2, 20).Value = 0.00008217591987769317

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cfh"
$ws.Cells.Item(3, 3).Value = "Sell"
$ws.Cells.Item(3, 4).Value = "M1"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.362097333333333
$ws.Cells.Item(3, 8).Value = 4.086292
$ws.Cells.Item(3, 9).Value = 0.0196292037450248
$ws.Cells.Item(3, 10).Value = 0.0196292037450248
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.974950666666667
$ws.Cells.Item(3, 14).Value = 8.924852
$ws.Cells.Item(3, 15).Value = 0.06056432486857137
$ws.Cells.Item(3, 16).Value = 0.06056432486857137
$ws.Cells.Item(3, 17).Value = 4.052172369864889
$ws.Cells.Item(3, 18).Value = 36.469551328784
$ws.Cells.Item(3, 19).Value = 0.00118882947252506
$ws.Cells.Item(3, 20).Value = 0.00118882947252506

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cfh"
$ws.Cells.Item(4, 3).Value = "Sell"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.362097333333333
$ws.Cells.Item(4, 8).Value = 4.086292
$ws.Cells.Item(4, 9).Value = 0.0196292037450248
$ws.Cells.Item(4, 10).Value = 0.0196292037450248
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 45.939923
$ws.Cells.Item(4, 14).Value = 137.819769
$ws.Cells.Item(4, 15).Value = 0.935249263856416
$ws.Cells.Item(4, 16).Value = 0.935249263856416
$ws.Cells.Item(4, 17).Value = 62.57464661183867
$ws.Cells.Item(4, 18).Value = 563.1718195065481
$ws.Cells.Item(4, 19).Value = 0.01835819835262205
$ws.Cells.Item(4, 20).Value = 0.01835819835262205

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cfh"
$ws.Cells.Item(5, 3).Value = "Sell"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 24.05951033333333
$ws.Cells.Item(5, 8).Value = 72.17853099999999
$ws.Cells.Item(5, 9).Value = 0.3467219403350491
$ws.Cells.Item(5, 10).Value = 0.346721940335049
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2056386666666667
$ws.Cells.Item(5, 14).Value = 0.616916
$ws.Cells.Item(5, 15).Value = 0.004186411275012692
$ws.Cells.Item(5, 16).Value = 0.004186411275012692
$ws.Cells.Item(5, 17).Value = 4.947565625599555
$ws.Cells.Item(5, 18).Value = 44.52809063039599
$ws.Cells.Item(5, 19).Value = 0.001451520640312927
$ws.Cells.Item(5, 20).Value = 0.001451520640312927

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cfh"
$ws.Cells.Item(6, 3).Value = "Sell"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 24.05951033333333
$ws.Cells.Item(6, 8).Value = 72.17853099999999
$ws.Cells.Item(6, 9).Value = 0.3467219403350491
$ws.Cells.Item(6, 10).Value = 0.346721940335049
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.974950666666667
$ws.Cells.Item(6, 14).Value = 8.924852
$ws.Cells.Item(6, 15).Value = 0.06056432486857137
$ws.Cells.Item(6, 16).Value = 0.06056432486857137
$ws.Cells.Item(6, 17).Value = 71.57585630582355
$ws.Cells.Item(6, 18).Value = 644.182706752412
$ws.Cells.Item(6, 19).Value = 0.02099898023351333
$ws.Cells.Item(6, 20).Value = 0.02099898023351333

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cfh"
$ws.Cells.Item(7, 3).Value = "Sell"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 24.05951033333333
$ws.Cells.Item(7, 8).Value = 72.17853099999999
$ws.Cells.Item(7, 9).Value = 0.3467219403350491
$ws.Cells.Item(7, 10).Value = 0.346721940335049
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 45.939923
$ws.Cells.Item(7, 14).Value = 137.819769
$ws.Cells.Item(7, 15).Value = 0.935249263856416
$ws.Cells.Item(7, 16).Value = 0.935249263856416
$ws.Cells.Item(7, 17).Value = 1105.292052131038
$ws.Cells.Item(7, 18).Value = 9947.628469179339
$ws.Cells.Item(7, 19).Value = 0.3242714394612228
$ws.Cells.Item(7, 20).Value = 0.3242714394612228

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Cfh"
$ws.Cells.Item(8, 3).Value = "Sell"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.711625
$ws.Cells.Item(8, 8).Value = 44.13487499999999
$ws.Cells.Item(8, 9).Value = 0.2120094338917045
$ws.Cells.Item(8, 10).Value = 0.2120094338917045
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.2056386666666667
$ws.Cells.Item(8, 14).Value = 0.616916
$ws.Cells.Item(8, 15).Value = 0.004186411275012692
$ws.Cells.Item(8, 16).Value = 0.004186411275012692
$ws.Cells.Item(8, 17).Value = 3.0252789495
$ws.Cells.Item(8, 18).Value = 27.2275105455
$ws.Cells.Item(8, 19).Value = 0.0008875586844532899
$ws.Cells.Item(8, 20).Value = 0.0008875586844532898

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Cfh"
$ws.Cells.Item(9, 3).Value = "Sell"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.711625
$ws.Cells.Item(9, 8).Value = 44.13487499999999
$ws.Cells.Item(9, 9).Value = 0.2120094338917045
$ws.Cells.Item(9, 10).Value = 0.2120094338917045
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.974950666666667
$ws.Cells.Item(9, 14).Value = 8.924852
$ws.Cells.Item(9, 15).Value = 0.06056432486857137
$ws.Cells.Item(9, 16).Value = 0.06056432486857137
$ws.Cells.Item(9, 17).Value = 43.7663586015
$ws.Cells.Item(9, 18).Value = 393.8972274134999
$ws.Cells.Item(9, 19).Value = 0.0128402082294191
$ws.Cells.Item(9, 20).Value = 0.0128402082294191

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Cfh"
$ws.Cells.Item(10, 3).Value = "Sell"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 14.711625
$ws.Cells.Item(10, 8).Value = 44.13487499999999
$ws.Cells.Item(10, 9).Value = 0.2120094338917045
$ws.Cells.Item(10, 10).Value = 0.2120094338917045
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 45.939923
$ws.Cells.Item(10, 14).Value = 137.819769
$ws.Cells.Item(10, 15).Value = 0.935249263856416
$ws.Cells.Item(10, 16).Value = 0.935249263856416
$ws.Cells.Item(10, 17).Value = 675.8509197048749
$ws.Cells.Item(10, 18).Value = 6082.658277343875
$ws.Cells.Item(10, 19).Value = 0.1982816669778322
$ws.Cells.Item(10, 20).Value = 0.1982816669778321

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Cfh"
$ws.Cells.Item(11, 3).Value = "Sell"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 27.616616
$ws.Cells.Item(11, 8).Value = 82.84984800000001
$ws.Cells.Item(11, 9).Value = 0.3979834399099074
$ws.Cells.Item(11, 10).Value = 0.3979834399099074
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.2056386666666667
$ws.Cells.Item(11, 14).Value = 0.616916
$ws.Cells.Item(11, 15).Value = 0.004186411275012692
$ws.Cells.Item(11, 16).Value = 0.004186411275012692
$ws.Cells.Item(11, 17).Value = 5.679044092085334
$ws.Cells.Item(11, 18).Value = 51.11139682876801
$ws.Cells.Item(11, 19).Value = 0.001666122360107173
$ws.Cells.Item(11, 20).Value = 0.001666122360107173

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Cfh"
$ws.Cells.Item(12, 3).Value = "Sell"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 27.616616
$ws.Cells.Item(12, 8).Value = 82.84984800000001
$ws.Cells.Item(12, 9).Value = 0.3979834399099074
$ws.Cells.Item(12, 10).Value = 0.3979834399099074
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.974950666666667
$ws.Cells.Item(12, 14).Value = 8.924852
$ws.Cells.Item(12, 15).Value = 0.06056432486857137
$ws.Cells.Item(12, 16).Value = 0.06056432486857137
$ws.Cells.Item(12, 17).Value = 82.15807018027733
$ws.Cells.Item(12, 18).Value = 739.422631622496
$ws.Cells.Item(12, 19).Value = 0.02410359834701518
$ws.Cells.Item(12, 20).Value = 0.02410359834701518

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Cfh"
$ws.Cells.Item(13, 3).Value = "Sell"
$ws.Cells.Item(13, 4).Value = "M2"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 27.616616
$ws.Cells.Item(13, 8).Value = 82.84984800000001
$ws.Cells.Item(13, 9).Value = 0.3979834399099074
$ws.Cells.Item(13, 10).Value = 0.3979834399099074
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 45.939923
$ws.Cells.Item(13, 14).Value = 137.819769
$ws.Cells.Item(13, 15).Value = 0.935249263856416
$ws.Cells.Item(13, 16).Value = 0.935249263856416
$ws.Cells.Item(13, 17).Value = 1268.705212560568
$ws.Cells.Item(13, 18).Value = 11418.34691304511
$ws.Cells.Item(13, 19).Value = 0.3722137192027851
$ws.Cells.Item(13, 20).Value = 0.3722137192027851

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Cfh"
$ws.Cells.Item(14, 3).Value = "Sell"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.641521
$ws.Cells.Item(14, 8).Value = 4.924563000000001
$ws.Cells.Item(14, 9).Value = 0.02365598211831425
$ws.Cells.Item(14, 10).Value = 0.02365598211831425
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.2056386666666667
$ws.Cells.Item(14, 14).Value = 0.616916
$ws.Cells.Item(14, 15).Value = 0.004186411275012692
$ws.Cells.Item(14, 16).Value = 0.004186411275012692
$ws.Cells.Item(14, 17).Value = 0.3375601897453334
$ws.Cells.Item(14, 18).Value = 3.038041707708001
$ws.Cells.Item(14, 19).Value = 0.0000990336702616094
$ws.Cells.Item(14, 20).Value = 0.0000990336702616094

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Cfh"
$ws.Cells.Item(15, 3).Value = "Sell"
$ws.Cells.Item(15, 4).Value = "M1"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.641521
$ws.Cells.Item(15, 8).Value = 4.924563000000001
$ws.Cells.Item(15, 9).Value = 0.02365598211831425
$ws.Cells.Item(15, 10).Value = 0.02365598211831425
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.974950666666667
$ws.Cells.Item(15, 14).Value = 8.924852
$ws.Cells.Item(15, 15).Value = 0.06056432486857137
$ws.Cells.Item(15, 16).Value = 0.06056432486857137
$ws.Cells.Item(15, 17).Value = 4.883443993297334
$ws.Cells.Item(15, 18).Value = 43.95099593967601
$ws.Cells.Item(15, 19).Value = 0.001432708586098699
$ws.Cells.Item(15, 20).Value = 0.001432708586098699

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Cfh"
$ws.Cells.Item(16, 3).Value = "Sell"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.641521
$ws.Cells.Item(16, 8).Value = 4.924563000000001
$ws.Cells.Item(16, 9).Value = 0.02365598211831425
$ws.Cells.Item(16, 10).Value = 0.02365598211831425
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 45.939923
$ws.Cells.Item(16, 14).Value = 137.819769
$ws.Cells.Item(16, 15).Value = 0.935249263856416
$ws.Cells.Item(16, 16).Value = 0.935249263856416
$ws.Cells.Item(16, 17).Value = 75.41134834288302
$ws.Cells.Item(16, 18).Value = 678.7021350859471
$ws.Cells.Item(16, 19).Value = 0.02212423986195394
$ws.Cells.Item(16, 20).Value = 0.02212423986195394
